$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-01-28 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-29 Monday", 2) | Out-Null

# Update each arithmetic expression cell in the table, addressed by row/col
$t = $d.Tables.Item(1)

$r = $t.Cell(1, 1).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "25-17="
$r = $t.Cell(1, 2).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "84-6="
$r = $t.Cell(1, 3).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "92-50="
$r = $t.Cell(1, 4).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "99-54="
$r = $t.Cell(1, 5).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "9+38="
$r = $t.Cell(2, 1).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "64-11="
$r = $t.Cell(2, 2).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "34+17="
$r = $t.Cell(2, 3).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "94-18="
$r = $t.Cell(2, 4).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "15-2="
$r = $t.Cell(2, 5).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "78-52="
$r = $t.Cell(3, 1).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "4+82="
$r = $t.Cell(3, 2).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "59-39="
$r = $t.Cell(3, 3).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "23-4="
$r = $t.Cell(3, 4).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "45+49="
$r = $t.Cell(3, 5).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "1+25="
$r = $t.Cell(4, 1).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "98-56="
$r = $t.Cell(4, 2).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "62-33="
$r = $t.Cell(4, 3).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "32+53="
$r = $t.Cell(4, 4).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "25+3="
$r = $t.Cell(4, 5).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "94-90="
$r = $t.Cell(5, 1).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "78+0="
$r = $t.Cell(5, 2).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "48+20="
$r = $t.Cell(5, 3).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "91-22="
$r = $t.Cell(5, 4).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "80-41="
$r = $t.Cell(5, 5).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "25-11="
$r = $t.Cell(6, 1).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "64-20="
$r = $t.Cell(6, 2).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "80-77="
$r = $t.Cell(6, 3).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "50+21="
$r = $t.Cell(6, 4).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "12+21="
$r = $t.Cell(6, 5).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "57-42="
$r = $t.Cell(7, 1).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "84-44="
$r = $t.Cell(7, 2).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "60-24="
$r = $t.Cell(7, 3).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "22+55="
$r = $t.Cell(7, 4).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "48-2="
$r = $t.Cell(7, 5).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "30+53="
$r = $t.Cell(8, 1).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "29+20="
$r = $t.Cell(8, 2).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "62-3="
$r = $t.Cell(8, 3).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "78-55="
$r = $t.Cell(8, 4).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "46-8="
$r = $t.Cell(8, 5).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "96-39="
$r = $t.Cell(9, 1).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "91-85="
$r = $t.Cell(9, 2).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "3-3="
$r = $t.Cell(9, 3).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "76+17="
$r = $t.Cell(9, 4).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "93+2="
$r = $t.Cell(9, 5).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "99-98="
$r = $t.Cell(10, 1).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "98-94="
$r = $t.Cell(10, 2).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "54-14="
$r = $t.Cell(10, 3).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "33+20="
$r = $t.Cell(10, 4).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "28+19="
$r = $t.Cell(10, 5).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "77-11="
$r = $t.Cell(11, 1).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "41+25="
$r = $t.Cell(11, 2).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "73-16="
$r = $t.Cell(11, 3).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "57+2="
$r = $t.Cell(11, 4).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "30+30="
$r = $t.Cell(11, 5).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "35+20="
$r = $t.Cell(12, 1).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "31+43="
$r = $t.Cell(12, 2).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "39+59="
$r = $t.Cell(12, 3).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "46-30="
$r = $t.Cell(12, 4).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "89+3="
$r = $t.Cell(12, 5).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "48+6="
$r = $t.Cell(13, 1).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "74-68="
$r = $t.Cell(13, 2).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "84-40="
$r = $t.Cell(13, 3).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "84-26="
$r = $t.Cell(13, 4).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "90-26="
$r = $t.Cell(13, 5).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "18+31="
$r = $t.Cell(14, 1).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "51+6="
$r = $t.Cell(14, 2).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "96-30="
$r = $t.Cell(14, 3).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "97-84="
$r = $t.Cell(14, 4).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "36+15="
$r = $t.Cell(14, 5).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "49+9="
$r = $t.Cell(15, 1).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "1+47="
$r = $t.Cell(15, 2).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "77+0="
$r = $t.Cell(15, 3).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "30+31="
$r = $t.Cell(15, 4).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "43+2="
$r = $t.Cell(15, 5).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "2+42="
$r = $t.Cell(16, 1).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "77+7="
$r = $t.Cell(16, 2).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "15+11="
$r = $t.Cell(16, 3).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "17+29="
$r = $t.Cell(16, 4).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "53-40="
$r = $t.Cell(16, 5).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "19+66="
$r = $t.Cell(17, 1).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "48-2="
$r = $t.Cell(17, 2).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "54+30="
$r = $t.Cell(17, 3).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "2+48="
$r = $t.Cell(17, 4).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "80-66="
$r = $t.Cell(17, 5).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "18-15="
$r = $t.Cell(18, 1).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "45+40="
$r = $t.Cell(18, 2).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "63-31="
$r = $t.Cell(18, 3).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "19+3="
$r = $t.Cell(18, 4).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "47+23="
$r = $t.Cell(18, 5).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "10+76="
$r = $t.Cell(19, 1).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "18+45="
$r = $t.Cell(19, 2).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "72+15="
$r = $t.Cell(19, 3).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "51-2="
$r = $t.Cell(19, 4).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "66-39="
$r = $t.Cell(19, 5).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "45-23="
$r = $t.Cell(20, 1).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "70-3="
$r = $t.Cell(20, 2).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "18+74="
$r = $t.Cell(20, 3).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "70+22="
$r = $t.Cell(20, 4).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "78+21="
$r = $t.Cell(20, 5).Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "62-19="
